# Update "upload table sql log aktivitas"
# - C2: 220927 -> 220220
# - E2: "ko gamasuk?" -> "masuk"
# - E3: "tapi ini bisa" -> "masuk"  (now a duplicate of E2, so the shared
#   string table collapses to a single "masuk" entry, matching the diff)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 220220
$ws.Range("E2").Value = "masuk"
$ws.Range("E3").Value = "masuk"
